$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.737.85"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.749.72"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5049"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").Value = "1.749.73"
$ws.Range("E11").Value = "  -5.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06913"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.77%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.464"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5949"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "25.752.19"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  +5.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006768"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.52%  "
$ws.Range("D22").Value = "1.974.16"
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.204"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.442"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.797"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "101.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  +7.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04474"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.29%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.659"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9871"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6018"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.683"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01543"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.929"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.41%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3770"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7356"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.40%  "
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05471"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1092"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.76%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.886"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.633"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
